$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names & links) ---
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Numeric-looking text cells (prices & percentages) ---
# Force Text format so Excel keeps these as literal strings rather than
# auto-converting them to numbers/percentages, then clear the formatting
# back off so no stray style is left on the cells.
$numCells = @("D2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D39","E39","D40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","E49","D50","E50","D51","E51")
foreach ($c in $numCells) { $ws.Range($c).NumberFormat = "@" }
$ws.Range("D2").Value = "304.00"
$ws.Range("D3").Value = "37.14"
$ws.Range("E3").Value = "3.37%"
$ws.Range("D4").Value = "5.036"
$ws.Range("E4").Value = "-0.87%"
$ws.Range("D5").Value = "0.07844"
$ws.Range("E5").Value = "0.24%"
$ws.Range("D6").Value = "2.204"
$ws.Range("E6").Value = "-3.74%"
$ws.Range("D7").Value = "8.006"
$ws.Range("E7").Value = "-0.52%"
$ws.Range("D8").Value = "4.028"
$ws.Range("E8").Value = "1.08%"
$ws.Range("D9").Value = "0.9301"
$ws.Range("E9").Value = "0.21%"
$ws.Range("D10").Value = "0.09909"
$ws.Range("E10").Value = "-1.52%"
$ws.Range("D11").Value = "0.1888"
$ws.Range("E11").Value = "3.04%"
$ws.Range("D12").Value = "0.08595"
$ws.Range("E12").Value = "0.05%"
$ws.Range("D13").Value = "0.03663"
$ws.Range("E13").Value = "7.81%"
$ws.Range("D14").Value = "0.09961"
$ws.Range("E14").Value = "0.61%"
$ws.Range("D15").Value = "0.001481"
$ws.Range("E15").Value = "0.34%"
$ws.Range("D16").Value = "0.005665"
$ws.Range("E16").Value = "-1.00%"
$ws.Range("D17").Value = "3.456"
$ws.Range("E17").Value = "-0.82%"
$ws.Range("D18").Value = "2.373"
$ws.Range("E18").Value = "17.33%"
$ws.Range("D19").Value = "0.3410"
$ws.Range("E19").Value = "-0.66%"
$ws.Range("D20").Value = "0.1309"
$ws.Range("E20").Value = "-1.26%"
$ws.Range("D21").Value = "4.768"
$ws.Range("E21").Value = "4.79%"
$ws.Range("D22").Value = "0.2205"
$ws.Range("E22").Value = "-0.66%"
$ws.Range("D23").Value = "0.04610"
$ws.Range("E23").Value = "-0.53%"
$ws.Range("D24").Value = "0.001256"
$ws.Range("E24").Value = "2.98%"
$ws.Range("D25").Value = "0.004481"
$ws.Range("E25").Value = "-0.08%"
$ws.Range("D26").Value = "0.0001403"
$ws.Range("E26").Value = "8.13%"
$ws.Range("D27").Value = "0.0002724"
$ws.Range("E27").Value = "-19.73%"
$ws.Range("D39").Value = "0.01849"
$ws.Range("E39").Value = "6.01%"
$ws.Range("D40").Value = "0.04777"
$ws.Range("D41").Value = "0.008046"
$ws.Range("E41").Value = "2.54%"
$ws.Range("D42").Value = "0.1404"
$ws.Range("E42").Value = "-0.77%"
$ws.Range("D43").Value = "0.007518"
$ws.Range("E43").Value = "-14.58%"
$ws.Range("D44").Value = "0.002132"
$ws.Range("E44").Value = "-3.78%"
$ws.Range("D45").Value = "0.01040"
$ws.Range("E45").Value = "13.45%"
$ws.Range("D46").Value = "0.00006301"
$ws.Range("E46").Value = "3.92%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.38%"
$ws.Range("D48").Value = "0.0005803"
$ws.Range("E48").Value = "0.04%"
$ws.Range("E49").Value = "551.18%"
$ws.Range("D50").Value = "0.002696"
$ws.Range("E50").Value = "0.35%"
$ws.Range("D51").Value = "0.00002105"
$ws.Range("E51").Value = "0.38%"
foreach ($c in $numCells) { $ws.Range($c).ClearFormats() }
